$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = "s"
